# Weekly update of the Betarraga (Vega Monumental Concepción) price sheet.
# Each existing dated record (rows 215-220) shifts down one slot to make
# room for the newest week, and the oldest record that falls off the
# bottom is re-appended as a new pair of rows (221-222) with its original
# values, extending the used range to A1:R222.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 221: Primera (copy of the data that used to live in row 219)
$ws.Range("A221").Value = 11
$ws.Range("B221").Value = "Vega Monumental Concepción"
$ws.Range("C221").Value = "Bíobío"
$ws.Range("D221").Value = 44273
$ws.Range("D221").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E221").Value = 8
$ws.Range("F221").Value = 100114014
$ws.Range("G221").Value = "Betarraga"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 600
$ws.Range("K221").Value = 600
$ws.Range("L221").Value = 700
$ws.Range("M221").Value = 650
$ws.Range("N221").Value = "`$/paquete 5 unidades"
$ws.Range("O221").Value = "Región Metropolitana"
$ws.Range("P221").Value = 130
$ws.Range("Q221").Value = 5
$ws.Range("R221").Value = "Hortaliza"

# New row 222: Segunda (copy of the data that used to live in row 220)
$ws.Range("A222").Value = 11
$ws.Range("B222").Value = "Vega Monumental Concepción"
$ws.Range("C222").Value = "Bíobío"
$ws.Range("D222").Value = 44273
$ws.Range("D222").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E222").Value = 8
$ws.Range("F222").Value = 100114014
$ws.Range("G222").Value = "Betarraga"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Segunda"
$ws.Range("J222").Value = 300
$ws.Range("K222").Value = 500
$ws.Range("L222").Value = 500
$ws.Range("M222").Value = 500
$ws.Range("N222").Value = "`$/paquete 5 unidades"
$ws.Range("O222").Value = "Región Metropolitana"
$ws.Range("P222").Value = 100
$ws.Range("Q222").Value = 5
$ws.Range("R222").Value = "Hortaliza"

# Row 220 (Segunda) -> becomes the old row 219/220 pattern shifted: date 44307, J 300->200
$ws.Range("D220").Value = 44307
$ws.Range("J220").Value = 200

# Row 219 (Primera) -> date 44307, J 600->500, M 650->660, P 130->132
$ws.Range("D219").Value = 44307
$ws.Range("J219").Value = 500
$ws.Range("M219").Value = 660
$ws.Range("P219").Value = 132

# Row 218 (Segunda) -> date 44383, J 200->300
$ws.Range("D218").Value = 44383
$ws.Range("J218").Value = 300

# Row 217 (Primera) -> date 44383, J 500->600, M 660->650, P 132->130
$ws.Range("D217").Value = 44383
$ws.Range("J217").Value = 600
$ws.Range("M217").Value = 650
$ws.Range("P217").Value = 130

# Row 216 (Segunda) -> date 44509 (values unchanged)
$ws.Range("D216").Value = 44509

# Row 215 (Primera) -> date 44509 (values unchanged)
$ws.Range("D215").Value = 44509
